$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add the new glossary term "framework " (note trailing space) as the next
# row beneath the existing single-column entries (row 32, column A).
$target = $ws.Range("A32")
$target.Value = "framework "

# Match the formatting used by the rest of the single-column glossary rows
# (wrap text, default font/fill) -- this reuses the existing style rather
# than creating a new one.
$target.WrapText = $true

# Keep the active selection on the newly added row, consistent with the
# authored workbook.
$target.Select()
